$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @('Backend Go Engineer', 'https://www.dice.com/job-detail/cc08a4b5-1ed3-4f51-af67-efec6373ac2b', 'New York, New York', 'Contract', '$73 - $83 per hour', 'Concent Software Solution LLC'),
    @('Backend Go Engineer', 'https://www.dice.com/job-detail/d88a759f-02b8-4ceb-bec3-12325eef67e8', 'New York, New York', 'Contract', '$65 - $75 per hour', 'Business Needs Inc.'),
    @('GoLang, Sr Engineer', 'https://www.dice.com/job-detail/467f0803-dae4-4cca-b5d3-3a3db36407c9', 'Remote', 'Contract', 'Depends on Experience', 'V-Soft Consulting Group, Inc')
)

$startRow = 58
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
